# Baker North 2015 site info: split "Date Sampled" into separate
# Month / Day / Year columns (massive MDY site template update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank columns before the old column E ("Date Sampled").
# This shifts Date Sampled .. Notes from E:N to H:Q and also shifts/updates
# the dependent dimension, dataValidation sqref, etc. automatically.
$ws.Columns("E:G").Insert()

# New header cells for the inserted columns.
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# New Month/Day/Year values for each data row, derived from the date
# that is now sitting in column H.
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 2015

$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 2015

$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 2015

$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 2015

# The new columns get a plain (non bestFit) custom width, close to the
# width used by the neighboring Transect column.
$ws.Columns("E:G").ColumnWidth = 8.33

# Selection moved from the old L8 to the newly inserted G5, and Excel no
# longer needs to scroll the view to column D.
$ws.Range("G5").Select() | Out-Null
